$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 12
$ws.Range("I8").Value = 12
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 36
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 103
$ws.Range("N8").ClearContents()
$ws.Range("H12").Value = 1964.1428
$ws.Range("J12").Value = 2434
$ws.Range("L12").Value = 2434
$ws.Range("N12").Value = -2774
$ws.Range("H62").Value = 3951.0527
$ws.Range("I62").Value = 2505
$ws.Range("K62").Value = 2505
$ws.Range("M62").Value = -1881
$ws.Range("H65").Value = 3951.0527
$ws.Range("I65").Value = 2505
$ws.Range("K65").Value = 12525
$ws.Range("M65").Value = -9405
$ws.Range("H96").Value = 2507.7144
$ws.Range("I96").Value = 2675.3333
$ws.Range("J96").Value = 2382
$ws.Range("K96").Value = 8025.999899999999
$ws.Range("L96").Value = 7146
$ws.Range("M96").Value = -6652.999899999999
$ws.Range("N96").Value = -9892
$ws.Range("H129").Value = 4585.143
$ws.Range("I129").Value = 1548
$ws.Range("K129").Value = 4644
$ws.Range("M129").Value = 356
$ws.Range("H132").Value = 2120.5881
$ws.Range("I132").Value = 2075.3572
$ws.Range("J132").Value = 2331.6667
$ws.Range("K132").Value = 6226.071599999999
$ws.Range("L132").Value = 6995.000100000001
$ws.Range("M132").Value = -3696.071599999999
$ws.Range("N132").Value = -12055.0001
$ws.Range("H137").Value = 2186.3225
$ws.Range("I137").Value = 1058.1177
$ws.Range("K137").Value = 3174.3531
$ws.Range("M137").Value = -624.3531000000003
$ws.Range("H138").Value = 6017.064
$ws.Range("I138").Value = 4367.2
$ws.Range("J138").Value = 6213.476
$ws.Range("K138").Value = 13101.6
$ws.Range("L138").Value = 18640.428
$ws.Range("M138").Value = -7961.599999999999
$ws.Range("N138").Value = -28920.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 834.4545000000001
$ws.Range("I2").Value = 734.6842
$ws.Range("K2").Value = 734.6842
$ws.Range("M2").Value = -621.6842
$ws.Range("H61").Value = 1381.2
$ws.Range("I61").Value = 1381.2
$ws.Range("K61").Value = 1381.2
$ws.Range("M61").Value = -1169.2
$ws.Range("H101").Value = 58999.5
$ws.Range("J101").Value = 58999.5
$ws.Range("L101").Value = 58999.5
$ws.Range("N101").Value = -65489.5
$ws.Range("H102").Value = 1962
$ws.Range("I102").Value = 1538.25
$ws.Range("J102").Value = 2640
$ws.Range("K102").Value = 1538.25
$ws.Range("L102").Value = 2640
$ws.Range("M102").Value = 83.75
$ws.Range("N102").Value = -5884
$ws.Range("H110").Value = 9347
$ws.Range("I110").Value = 11463
$ws.Range("J110").Value = 2999
$ws.Range("K110").Value = 11463
$ws.Range("L110").Value = 2999
$ws.Range("M110").Value = -9418
$ws.Range("N110").Value = -7089
$ws.Range("H116").Value = 834.4545000000001
$ws.Range("I116").Value = 734.6842
$ws.Range("K116").Value = 734.6842
$ws.Range("M116").Value = 1559.3158
$ws.Range("H132").Value = 2294.257
$ws.Range("I132").Value = 2142.5312
$ws.Range("J132").Value = 3912.6667
$ws.Range("K132").Value = 6427.5936
$ws.Range("L132").Value = 11738.0001
$ws.Range("M132").Value = -3897.5936
$ws.Range("N132").Value = -16798.0001
$ws.Range("H136").Value = 1381.2
$ws.Range("I136").Value = 1381.2
$ws.Range("K136").Value = 4143.6
$ws.Range("M136").Value = -1593.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 834.4545000000001
$ws.Range("I3").Value = 734.6842
$ws.Range("K3").Value = 734.6842
$ws.Range("M3").Value = -620.6842
$ws.Range("H49").Value = 25000
$ws.Range("J49").Value = 25000
$ws.Range("L49").Value = 25000
$ws.Range("N49").Value = -25478
$ws.Range("H64").Value = 1034.5385
$ws.Range("I64").Value = 1029.2
$ws.Range("J64").Value = 1037.875
$ws.Range("K64").Value = 1029.2
$ws.Range("L64").Value = 1037.875
$ws.Range("M64").Value = -804.2
$ws.Range("N64").Value = -1487.875
$ws.Range("H67").Value = 1034.5385
$ws.Range("I67").Value = 1029.2
$ws.Range("J67").Value = 1037.875
$ws.Range("K67").Value = 1029.2
$ws.Range("L67").Value = 1037.875
$ws.Range("M67").Value = -249.2
$ws.Range("N67").Value = -2597.875
$ws.Range("H134").Value = 1430.4783
$ws.Range("I134").Value = 770.6
$ws.Range("K134").Value = 2311.8
$ws.Range("M134").Value = 223.1999999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 272.1111
$ws.Range("I22").Value = 224.83333
$ws.Range("K22").Value = 224.83333
$ws.Range("M22").Value = 125.16667
$ws.Range("H134").Value = 2032.1765
$ws.Range("I134").Value = 1253.3334
$ws.Range("J134").Value = 5036.2856
$ws.Range("K134").Value = 3760.0002
$ws.Range("L134").Value = 15108.8568
$ws.Range("M134").Value = -1225.0002
$ws.Range("N134").Value = -20178.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 641
$ws.Range("I5").Value = 498.66666
$ws.Range("K5").Value = 1495.99998
$ws.Range("M5").Value = -1383.99998
$ws.Range("H33").Value = 780.5714
$ws.Range("I33").Value = 1758.6666
$ws.Range("J33").Value = 47
$ws.Range("K33").Value = 10551.9996
$ws.Range("L33").Value = 282
$ws.Range("M33").Value = -10268.9996
$ws.Range("N33").Value = -848
$ws.Range("H34").Value = 1949.2307
$ws.Range("I34").Value = 1207.4286
$ws.Range("J34").Value = 2814.6667
$ws.Range("K34").Value = 3622.2858
$ws.Range("L34").Value = 8444.000100000001
$ws.Range("M34").Value = -3538.2858
$ws.Range("N34").Value = -8612.000100000001
$ws.Range("H122").Value = 933.7778
$ws.Range("I122").Value = 771.6
$ws.Range("K122").Value = 6944.400000000001
$ws.Range("M122").Value = -4494.400000000001
$ws.Range("H131").Value = 2949.8667
$ws.Range("J131").Value = 3408.4546
$ws.Range("L131").Value = 10225.3638
$ws.Range("N131").Value = -20305.3638
$ws.Range("H132").Value = 3103.6365
$ws.Range("I132").Value = 956.8333
$ws.Range("K132").Value = 8611.4997
$ws.Range("M132").Value = -6081.4997
$ws.Range("H135").Value = 641
$ws.Range("I135").Value = 498.66666
$ws.Range("K135").Value = 4487.99994
$ws.Range("M135").Value = -1952.99994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5663.385
$ws.Range("I80").Value = 3626.111
$ws.Range("K80").Value = 3626.111
$ws.Range("M80").Value = -2628.111
$ws.Range("H83").Value = 5663.385
$ws.Range("I83").Value = 3626.111
$ws.Range("K83").Value = 18130.555
$ws.Range("M83").Value = -13138.555
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 480896.8
$ws.Range("I122").Value = 60508.65
$ws.Range("J122").Value = 1671996.6
$ws.Range("K122").Value = 181525.95
$ws.Range("L122").Value = 5015989.800000001
$ws.Range("M122").Value = -179075.95
$ws.Range("N122").Value = -5020889.800000001
$ws.Range("H126").Value = 2609.158
$ws.Range("I126").Value = 1566.25
$ws.Range("K126").Value = 4698.75
$ws.Range("M126").Value = -2228.75
$ws.Range("H132").Value = 4655.1113
$ws.Range("I132").Value = 3644.7144
$ws.Range("J132").Value = 8191.5
$ws.Range("K132").Value = 10934.1432
$ws.Range("L132").Value = 24574.5
$ws.Range("M132").Value = -8404.143199999999
$ws.Range("N132").Value = -29634.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 40019
$ws.Range("J28").Value = 40019
$ws.Range("L28").Value = 40019
$ws.Range("N28").Value = -40715

